$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add older catch limits -- fill in landings (B) and target TAC (D) values
# that were previously blank, and turn the "Redfish - 500s" placeholder
# row into an actual "Redfish" data row.

$ws.Range("B13").Value = 8501
$ws.Range("D13").Value = 56673.333333333336

$ws.Range("A15").Value = "Redfish"
$ws.Range("B15").Value = 723
$ws.Range("D15").Value = 72300

$ws.Range("D17").Value = 57766.666666666672

$ws.Range("D19").Value = 55272.222222222226

$ws.Range("D25").Value = 54009.090909090912

# Cursor/selection moved to K12 in the saved file
$ws.Range("K12").Select()
